$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")
$ws.Activate()

# Order matters for shared-string slot re-use: first free up the old "FAIL" string
# by changing the row 7 Results cell (its sole reference) to the existing "SKIP" text,
# then claim that freed slot with the new "PASS" text on row 12, and finally append the
# brand new "N" text used by most of the Runmode (column D) cells.
$ws.Cells.Item(7, 5).Value = "SKIP"
$ws.Cells.Item(12, 5).Value = "PASS"

# Column D ("Runmode") switches from "Y" to "N" for every row except 12 and 15,
# which keep running ("Y").
for ($r = 2; $r -le 20; $r++) {
    if ($r -ne 12 -and $r -ne 15) {
        $ws.Cells.Item($r, 4).Value = "N"
    }
}

# Update the sheet view: scroll so row 8 is at the top and select B12.
$win = $excel.ActiveWindow
$win.ScrollRow = 8
$win.ScrollColumn = 1
$ws.Range("B12").Select()
